# Adds an "Active" (Y) column to each of the 4 sheets of the network-flow
# configuration workbook: Arcs (Table1), Node_start, Node_end, Nodes.
# Also bumps one data value on Nodes (D2: 1.65 -> 1.75) and refreshes
# selections / filter-database defined names to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Arcs sheet: genuine Excel Table (Table1). Add a 6th column "Active"
#    and fill every data row with "Y", matching the existing banded
#    (odd/even) row style used by the rest of the table.
# ---------------------------------------------------------------------
$wsArcs = $wb.Worksheets.Item("Arcs")
$tbl = $wsArcs.ListObjects.Item("Table1")
$tbl.ListColumns.Add() | Out-Null
$wsArcs.Range("F1").Value2 = "Active"
for ($r = 2; $r -le 11; $r++) {
    $wsArcs.Cells.Item($r, 1).Copy($wsArcs.Cells.Item($r, 6))
    $wsArcs.Cells.Item($r, 6).Value2 = "Y"
}

# ---------------------------------------------------------------------
# 2) Node_start sheet: plain AutoFilter range A1:E2 -> A1:F2.
# ---------------------------------------------------------------------
$wsStart = $wb.Worksheets.Item("Node_start")
$wsStart.Range("F1").Value2 = "Active"
$wsStart.Range("F2").Value2 = "Y"
$wsStart.Range("A1:F2").AutoFilter() | Out-Null
$wsStart.Range("E2").Select()

# ---------------------------------------------------------------------
# 3) Node_end sheet: plain AutoFilter range A1:D2 -> A1:E2.
# ---------------------------------------------------------------------
$wsEnd = $wb.Worksheets.Item("Node_end")
$wsEnd.Range("E1").Value2 = "Active"
$wsEnd.Range("E2").Value2 = "Y"
$wsEnd.Range("A1:E2").AutoFilter() | Out-Null
$wsEnd.Range("D2").Select()

# ---------------------------------------------------------------------
# 4) Nodes sheet: plain AutoFilter range A1:D7 -> A1:E7, plus a data fix
#    on D2 (1.65 -> 1.75).
# ---------------------------------------------------------------------
$wsNodes = $wb.Worksheets.Item("Nodes")
$wsNodes.Range("D2").Value2 = 1.75
$wsNodes.Range("E1").Value2 = "Active"
for ($r = 2; $r -le 7; $r++) {
    $wsNodes.Cells.Item($r, 5).Value2 = "Y"
}
$wsNodes.Range("A1:E7").AutoFilter() | Out-Null
$wsNodes.Range("E3").Select()

# ---------------------------------------------------------------------
# 5) Refresh the stale workbook-level _FilterDatabase defined names so
#    they track the new ranges (Arcs' is left untouched, matching the
#    observed behaviour of the authored edit).
# ---------------------------------------------------------------------
$wb.Names.Item("Node_end!_FilterDatabase").RefersTo = "=Node_end!`$A`$1:`$E`$2"
$wb.Names.Item("Node_start!_FilterDatabase").RefersTo = "=Node_start!`$A`$1:`$F`$2"
$wb.Names.Item("Nodes!_FilterDatabase").RefersTo = "=Nodes!`$A`$1:`$E`$7"

# ---------------------------------------------------------------------
# 6) Leave Arcs as the active sheet/selection, as in the source edit
#    (sheet1 keeps tabSelected="1", selection moves from M8 to E4).
# ---------------------------------------------------------------------
$wsArcs.Activate()
$wsArcs.Range("E4").Select()

Write-Output "edit complete"
